# Updates the EC (Estado de Cuenta) worker/period table so that rows are
# grouped by period instead of by worker:
#   Row 16 (period 2001 / Victor)  - keep worker, change period 2002 -> 2001
#   Row 17 (period 2001 / Oralys)  - becomes Oralys / 45371566 / 2001
#   Row 18 (period 2002 / Victor)  - becomes Victor / 1047382169 / 2002
#   Row 19 (period 2002 / Oralys)  - keep worker, change period 2001 -> 2002

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: only the "Periodo Mora" changes (2002 -> 2001)
$ws.Range("E16").Value = "2001"

# Row 17: Oralys Barrios Bello, doc 45371566, periodo 2001
$ws.Range("C17").Value = "45371566"
$ws.Range("D17").Value = "ORALYS BARRIOS BELLO"
$ws.Range("E17").Value = "2001"

# Row 18: Victor Augusto Padilla Ballestas, doc 1047382169, periodo 2002
$ws.Range("C18").Value = "1047382169"
$ws.Range("D18").Value = "VICTOR AUGUSTO PADILLA BALLESTAS"
$ws.Range("E18").Value = "2002"

# Row 19: only the "Periodo Mora" changes (2001 -> 2002)
$ws.Range("E19").Value = "2002"
